$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $val
    $range.Style = "Normal"
}

Set-TextCell "D2" "63.283.09"
Set-TextCell "E2" "  +0.80%  "
Set-TextCell "D3" "3.400.28"
Set-TextCell "E3" "  +1.55%  "
Set-TextCell "E4" "  +0.02%  "
Set-TextCell "D5" "567.22"
Set-TextCell "E5" "  +0.82%  "
Set-TextCell "D6" "155.68"
Set-TextCell "E6" "  +2.26%  "
Set-TextCell "D7" "0.999"
Set-TextCell "E7" "  -0.05%  "
Set-TextCell "D8" "3.398.10"
Set-TextCell "E8" "  +1.32%  "
Set-TextCell "E9" "  +2.42%  "
Set-TextCell "E10" "  -0.86%  "
Set-TextCell "E11" "  +3.24%  "
Set-TextCell "E12" "  -0.75%  "
Set-TextCell "D13" "3.984.90"
Set-TextCell "E13" "  +1.59%  "
Set-TextCell "E14" "  -3.04%  "
Set-TextCell "D15" "0.0000192"
Set-TextCell "E15" "  +7.07%  "
Set-TextCell "D16" "27.16"
Set-TextCell "E16" "  +0.83%  "
Set-TextCell "D17" "63.347.53"
Set-TextCell "E17" "  +0.93%  "
Set-TextCell "D18" "3.391.37"
Set-TextCell "E18" "  +2.21%  "
Set-TextCell "E19" "  -1.78%  "
Set-TextCell "D20" "14.05"
Set-TextCell "E20" "  +1.46%  "
Set-TextCell "D21" "378.81"
Set-TextCell "E21" "  -1.58%  "
Set-TextCell "D22" "8.05"
Set-TextCell "E22" "  -3.83%  "
Set-TextCell "D23" "1.00"
Set-TextCell "E23" "  +0.00%  "
Set-TextCell "D24" "71.43"
Set-TextCell "E24" "  +1.64%  "
Set-TextCell "D25" "0.528"
Set-TextCell "E25" "  -1.56%  "
Set-TextCell "E26" "  +25.77%  "
Set-TextCell "D27" "9.44"
Set-TextCell "E27" "  +5.76%  "
Set-TextCell "D28" "0.179"
Set-TextCell "E28" "  +0.03%  "
Set-TextCell "D30" "6.03"
Set-TextCell "E30" "  +7.86%  "
Set-TextCell "E31" "  +3.98%  "
Set-TextCell "E32" "  +0.62%  "
Set-TextCell "B33" "RenderToken"
Set-TextCell "C33" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D33" "6.41"
Set-TextCell "E33" "  -2.80%  "
Set-TextCell "B34" "EthereumClassic"
Set-TextCell "C34" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D34" "23.16"
Set-TextCell "E34" "  +0.99%  "
Set-TextCell "D35" "0.999"
Set-TextCell "E35" "  +0.07%  "
Set-TextCell "D36" "6.78"
Set-TextCell "E36" "  +1.30%  "
Set-TextCell "D37" "159.66"
Set-TextCell "E37" "  -0.21%  "
Set-TextCell "E38" "  -1.90%  "
Set-TextCell "D39" "2.956.77"
Set-TextCell "E39" "  +4.69%  "
Set-TextCell "D40" "26.94"
Set-TextCell "E40" "  +0.01%  "
Set-TextCell "B41" "Hedera"
Set-TextCell "C41" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D41" "0.0752"
Set-TextCell "E41" "  +1.47%  "
Set-TextCell "B42" "Stacks"
Set-TextCell "C42" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell "D42" "1.82"
Set-TextCell "E42" "  -3.95%  "
Set-TextCell "D43" "0.0315"
Set-TextCell "E43" "  +0.82%  "
Set-TextCell "D44" "41.73"
Set-TextCell "E44" "  +2.98%  "
Set-TextCell "D45" "0.761"
Set-TextCell "E45" "  +2.12%  "
Set-TextCell "D46" "4.32"
Set-TextCell "E46" "  +1.37%  "
Set-TextCell "D47" "23.25"
Set-TextCell "E47" "  +6.00%  "
Set-TextCell "E48" "  +3.25%  "
Set-TextCell "D49" "2.21"
Set-TextCell "E49" "  +22.94%  "
Set-TextCell "D50" "0.837"
Set-TextCell "E50" "  +4.78%  "
Set-TextCell "D51" "6.34"
Set-TextCell "E51" "  +0.73%  "
